# Update the "Creation_date" column (C) values to reflect the passage of
# one additional month for the proposals that are not yet a year old
# (rows pulled from the non-proposal / proposal links refresh).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C5"  = "Created 4 months ago"
    "C6"  = "Created 4 months ago"
    "C8"  = "Created 5 months ago"
    "C9"  = "Created 5 months ago"
    "C11" = "Created 6 months ago"
    "C13" = "Created 9 months ago"
    "C15" = "Created 10 months ago"
    "C17" = "Created a year ago"
    "C18" = "Created a year ago"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
